# Linear Model Results.xlsx - apply "implemented 3d resnet and have written
# down linear and nonlinear model results" edit:
#   * insert a new "Model Details" column between Outputs and Input Details
#   * append two new result rows (3D ResNet model) at the bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column C ("Model Details"), shifting the old C..H data
#    (Input Details, Training Loss, Testing Loss, Training Accuracy,
#    Testing Accuracy, Notes) one column to the right, into D..I.
# ---------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Insert() | Out-Null

$ws.Range("C1").Value = "Model Details"

# Give the new column the same width as the "Outputs" column next to it.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# ---------------------------------------------------------------------
# 2. Append the two new rows of results (3D ResNet runs) below the
#    existing table, which now ends at row 9.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "4 word features"
$ws.Range("B10").Value = "4 word features summed"
$ws.Range("C10").Value = "Fully connect layer, no activations"
$ws.Range("D10").Value = "^"
$ws.Range("E10").Value = 2869
$ws.Range("F10").Value = 1544
$ws.Range("G10").Value = 0.64
$ws.Range("H10").Value = 0.64
$ws.Range("I10").Value = "1 fold, 1000 epochs, LR=1e-7, momentum=0.9"

$ws.Range("A11").Value = "^"
$ws.Range("B11").Value = "^"
$ws.Range("C11").Value = "^"
$ws.Range("E11").Value = 1181
$ws.Range("F11").Value = 1098
$ws.Range("G11").Value = 0.88
$ws.Range("H11").Value = 0.79
$ws.Range("I11").Value = "1 fold, 1000 epochs, LR=1e-3, momentum=0.9"

# ---------------------------------------------------------------------
# 3. Leave the selection where the author ended up after typing the
#    final row (one cell past the new table, column I / row 12).
# ---------------------------------------------------------------------
$ws.Range("I12").Select() | Out-Null
